$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing odds values that changed between snapshots
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.62
$ws.Range("G5").Value = 1.4
$ws.Range("I5").Value = 8
$ws.Range("K5").Value = 2.38
$ws.Range("Z5").Value = 9
$ws.Range("AJ5").Value = 23
$ws.Range("G6").Value = 1.91
$ws.Range("H6").Value = 3.7
$ws.Range("I6").Value = 3.9
$ws.Range("K6").Value = 2.2
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.75
$ws.Range("AA6").Value = 15
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 10
$ws.Range("AE6").Value = 15
$ws.Range("AH6").Value = 11
$ws.Range("AO6").Value = 10
$ws.Range("AT6").Value = 2.75
$ws.Range("AW6").Value = 6
$ws.Range("AZ6").Value = 67
$ws.Range("BB6").Value = 201
$ws.Range("W7").Value = 8.5
$ws.Range("AH7").Value = 9.5
$ws.Range("AK7").Value = 29
$ws.Range("AM7").Value = 29
$ws.Range("G8").Value = 2.1
$ws.Range("I8").Value = 3.4
$ws.Range("J8").Value = 2.88
$ws.Range("L8").Value = 4.33
$ws.Range("N8").Value = 8.5
$ws.Range("W8").Value = 6.5
$ws.Range("Y8").Value = 9
$ws.Range("Z8").Value = 19
$ws.Range("AH8").Value = 9
$ws.Range("AI8").Value = 17
$ws.Range("AK8").Value = 41
$ws.Range("AO8").Value = 12
$ws.Range("AV8").Value = 67
$ws.Range("AW8").Value = 5.5
$ws.Range("AX8").Value = 21
$ws.Range("G10").Value = 2.8
$ws.Range("I10").Value = 2.3
$ws.Range("J10").Value = 3.4
$ws.Range("AD10").Value = 6.5
$ws.Range("AK10").Value = 23
$ws.Range("AL10").Value = 19
$ws.Range("N16").Value = 9
$ws.Range("Q18").Value = 1.95
$ws.Range("R18").Value = 1.9

# Add new row 21 (Wanderers vs Miramar, Uruguay Primera Division)
$ws.Range("A21").Value = "UeSMa6cR"
$ws.Range("B21").Value = "26/11/2024"
$ws.Range("C21").Value = "18:45"
$ws.Range("D21").Value = "URUGUAY - PRIMERA DIVISION"
$ws.Range("E21").Value = "Wanderers"
$ws.Range("F21").Value = "Miramar"
$ws.Range("G21").Value = 2.25
$ws.Range("H21").Value = 3.4
$ws.Range("I21").Value = 3.1
$ws.Range("J21").Value = 3
$ws.Range("K21").Value = 2.2
$ws.Range("L21").Value = 3.6
$ws.Range("M21").Value = 1.05
$ws.Range("N21").Value = 11
$ws.Range("O21").Value = 1.29
$ws.Range("P21").Value = 3.5
$ws.Range("Q21").Value = 1.98
$ws.Range("R21").Value = 1.88
$ws.Range("S21").Value = 1.4
$ws.Range("T21").Value = 2.75
$ws.Range("U21").Value = 1.73
$ws.Range("V21").Value = 2
$ws.Range("W21").Value = 8.5
$ws.Range("X21").Value = 11
$ws.Range("Y21").Value = 9.5
$ws.Range("Z21").Value = 21
$ws.Range("AA21").Value = 19
$ws.Range("AB21").Value = 26
$ws.Range("AC21").Value = 11
$ws.Range("AD21").Value = 6.5
$ws.Range("AE21").Value = 13
$ws.Range("AF21").Value = 41
$ws.Range("AG21").Value = 201
$ws.Range("AH21").Value = 10
$ws.Range("AI21").Value = 15
$ws.Range("AJ21").Value = 11
$ws.Range("AK21").Value = 34
$ws.Range("AL21").Value = 23
$ws.Range("AM21").Value = 34
$ws.Range("AN21").Value = 4.33
$ws.Range("AO21").Value = 13
$ws.Range("AP21").Value = 21
$ws.Range("AQ21").Value = 41
$ws.Range("AR21").Value = 51
$ws.Range("AS21").Value = 151
$ws.Range("AT21").Value = 2.75
$ws.Range("AU21").Value = 8
$ws.Range("AV21").Value = 51
$ws.Range("AW21").Value = 5
$ws.Range("AX21").Value = 17
$ws.Range("AY21").Value = 26
$ws.Range("AZ21").Value = 51
$ws.Range("BA21").Value = 67
$ws.Range("BB21").Value = 151
$ws.Range("BC21").Value = 51
$ws.Range("BD21").Value = 51
